# Fix: Elective lecture and tutorial scheduling
# Update classroom codes in the "Lecture Slot - Classroom" (D) and
# "Tutorial Slot - Classroom" (E) columns of the elective baskets table
# (rows 20-35) on all three timetable sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Regular_Timetable", "PreMid_Timetable", "PostMid_Timetable")

# Row number -> [ new D value, new E value ]
$updates = @{
    20 = @("Mon 09:00-10:30 [C101], Wed 13:00-14:30 [C101]", "Tue 14:30-15:30 [C101]")
    21 = @("Mon 09:00-10:30 [C102], Wed 13:00-14:30 [C102]", "Tue 14:30-15:30 [C102]")
    22 = @("Mon 09:00-10:30 [C104], Wed 13:00-14:30 [C104]", "Tue 14:30-15:30 [C104]")
    23 = @("Tue 09:00-10:30 [C101], Thu 13:00-14:30 [C101]", "Wed 14:30-15:30 [C101]")
    24 = @("Tue 09:00-10:30 [C102], Thu 13:00-14:30 [C102]", "Wed 14:30-15:30 [C102]")
    25 = @("Tue 09:00-10:30 [C104], Thu 13:00-14:30 [C104]", "Wed 14:30-15:30 [C104]")
    26 = @("Tue 09:00-10:30 [C202], Thu 13:00-14:30 [C202]", "Wed 14:30-15:30 [C202]")
    27 = @("Mon 10:30-12:00 [C101], Wed 10:30-12:00 [C101]", "Thu 14:30-15:30 [C101]")
    28 = @("Mon 10:30-12:00 [C102], Wed 10:30-12:00 [C102]", "Thu 14:30-15:30 [C102]")
    29 = @("Mon 10:30-12:00 [C104], Wed 10:30-12:00 [C104]", "Thu 14:30-15:30 [C104]")
    30 = @("Mon 10:30-12:00 [C202], Wed 10:30-12:00 [C202]", "Thu 14:30-15:30 [C202]")
    31 = @("Mon 10:30-12:00 [C203], Wed 10:30-12:00 [C203]", "Thu 14:30-15:30 [C203]")
    32 = @("Tue 15:30-17:00 [C101], Thu 15:30-17:00 [C101]", "Fri 14:30-15:30 [C101]")
    33 = @("Tue 15:30-17:00 [C102], Thu 15:30-17:00 [C102]", "Fri 14:30-15:30 [C102]")
    34 = @("Tue 15:30-17:00 [C104], Thu 15:30-17:00 [C104]", "Fri 14:30-15:30 [C104]")
    35 = @("Tue 15:30-17:00 [C202], Thu 15:30-17:00 [C202]", "Fri 14:30-15:30 [C202]")
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $vals = $updates[$row]
        $ws.Range("D$row").Value = $vals[0]
        $ws.Range("E$row").Value = $vals[1]
    }
}
